$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of tracker data (row 3)
$ws.Range("A3").Value = "G2"
$ws.Range("B3").Value = "Test1"
$ws.Range("C3").Value = "Daily"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 45860
$ws.Range("F3").Value = 30

# Match the date formatting used on the existing DateAdded column (E2)
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
